# Generate Report for Handoff
# File "6603f24f-90f7-44ae-8799-dd3d11673c6d.md" moves from "In Translation"
# to "Ready for handoff" status, with updated handoff timestamps.

$wb = $excel.ActiveWorkbook

# --- Overview sheet (row 4 = 6603f24f-90f7-44ae-8799-dd3d11673c6d.md) ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B4").Value = "Ready for handoff"
$overview.Range("C4").Value = "Ready for handoff"
$overview.Range("D4").Value = "2016-03-23 04:01:26"

# --- zh-cn sheet (row 4 = 6603f24f-90f7-44ae-8799-dd3d11673c6d.md) ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C4").Value = "Ready for handoff"
$zhcn.Range("E4").Value = "2016-03-23 04:01:08"

# --- de-de sheet (row 4 = 6603f24f-90f7-44ae-8799-dd3d11673c6d.md) ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C4").Value = "Ready for handoff"
$dede.Range("E4").Value = "2016-03-23 04:01:26"
